{"js": "// Replace the arithmetic-problem text runs throughout the document body.\n// Each entry is [oldText, newText]; all pairs are unique on both sides so a\n// simple search-and-replace per pair is safe (no collisions/re-matches).\nconst replacements = [\n  [\"268\u00d77=\", \"334\u00d72=\"],\n  [\"571\u00d74=\", \"905\u00d74=\"],\n  [\"659\u00d72=\", \"342\u00d78=\"],\n  [\"285\u00d76=\", \"226\u00d76=\"],\n  [\"575\u00d78=\", \"879\u00d79=\"],\n  [\"841\u00d73=\", \"699\u00d79=\"],\n  [\"677\u00d76=\", \"606\u00d77=\"],\n  [\"628\u00d76=\", \"903\u00d73=\"],\n  [\"974\u00d74=\", \"375\u00d76=\"],\n  [\"190\u00d78=\", \"259\u00d76=\"],\n  [\"405\u00d74=\", \"857\u00d76=\"],\n  [\"991\u00d74=\", \"136\u00d77=\"],\n  [\"476\u00d79=\", \"568\u00d76=\"],\n  [\"229\u00d74=\", \"263\u00d75=\"],\n  [\"837\u00d77=\", \"384\u00d72=\"],\n  [\"440\u00d76=\", \"473\u00d73=\"],\n  [\"749\u00d73=\", \"774\u00d74=\"],\n  [\"131\u00d74=\", \"953\u00d79=\"],\n  [\"264\u00d77=\", \"775\u00d75=\"],\n  [\"254\u00d79=\", \"687\u00d76=\"],\n  [\"981\u00d74=\", \"841\u00d72=\"],\n  [\"438\u00d79=\", \"814\u00d73=\"],\n  [\"277\u00d73=\", \"704\u00d74=\"],\n  [\"685\u00d72=\", \"480\u00d72=\"],\n  [\"996\u00d77=\", \"612\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the arithmetic-problem text runs throughout the document body.\n# Each pair is unique on both sides, so a single Find/Replace (wdReplaceAll)\n# per pair safely retargets exactly one run without touching anything else.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"268\u00d77=\", \"334\u00d72=\"),\n    @(\"571\u00d74=\", \"905\u00d74=\"),\n    @(\"659\u00d72=\", \"342\u00d78=\"),\n    @(\"285\u00d76=\", \"226\u00d76=\"),\n    @(\"575\u00d78=\", \"879\u00d79=\"),\n    @(\"841\u00d73=\", \"699\u00d79=\"),\n    @(\"677\u00d76=\", \"606\u00d77=\"),\n    @(\"628\u00d76=\", \"903\u00d73=\"),\n    @(\"974\u00d74=\", \"375\u00d76=\"),\n    @(\"190\u00d78=\", \"259\u00d76=\"),\n    @(\"405\u00d74=\", \"857\u00d76=\"),\n    @(\"991\u00d74=\", \"136\u00d77=\"),\n    @(\"476\u00d79=\", \"568\u00d76=\"),\n    @(\"229\u00d74=\", \"263\u00d75=\"),\n    @(\"837\u00d77=\", \"384\u00d72=\"),\n    @(\"440\u00d76=\", \"473\u00d73=\"),\n    @(\"749\u00d73=\", \"774\u00d74=\"),\n    @(\"131\u00d74=\", \"953\u00d79=\"),\n    @(\"264\u00d77=\", \"775\u00d75=\"),\n    @(\"254\u00d79=\", \"687\u00d76=\"),\n    @(\"981\u00d74=\", \"841\u00d72=\"),\n    @(\"438\u00d79=\", \"814\u00d73=\"),\n    @(\"277\u00d73=\", \"704\u00d74=\"),\n    @(\"685\u00d72=\", \"480\u00d72=\"),\n    @(\"996\u00d77=\", \"612\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
